$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hyperlinks (F2: emp1@xyz.com, F4: emp1@abc.com)
$ws.Hyperlinks.Delete()

# Clear the hyperlink cell formatting (remove the "Hyperlink" style reference)
$ws.Range("F2").ClearFormats()
$ws.Range("F4").ClearFormats()

# Update Portfolio Company values (column D): XYZ -> Apple, ABC -> MSFT
$ws.Range("D2").Value = "Apple"
$ws.Range("D3").Value = "Apple"
$ws.Range("D4").Value = "MSFT"
$ws.Range("D5").Value = "MSFT"

# Delete columns E (Pan) and F (Primary Email *); old G/H (Category*/Sub Category*) shift to E/F
$ws.Columns("E:F").Delete()

# Update selection
$ws.Range("D6").Select()

# Remove the now-unused "Hyperlink" cell style definition
$wb.Styles("Hyperlink").Delete()

$wb.Save()
